# Apply the cryptos-list price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write $text into $cellRef as a literal TEXT value (never let Excel
# auto-convert number-looking strings like "301.06" into a numeric cell).
# We build the text in an off-sheet scratch cell via a formula that evaluates
# to a string (="...") and then copy/paste-special VALUES ONLY into the real
# destination cell. Pasting a value that is already typed as Text by the engine
# keeps the destination cell typed as Text too, instead of re-parsing the raw
# characters as a number the way a direct .Value assignment would.
function Set-TextValue($cellRef, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.ClearContents()
}

$ws.Range("D2").Value = "43.112.75"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "2.306.62"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "301.06"
$ws.Range("E5").Value = "  -0.40%  "
Set-TextValue "D6" "97.81"
$ws.Range("E6").Value = "  -2.35%  "
Set-TextValue "D7" "0.520"
$ws.Range("E7").Value = "  +2.92%  "
$ws.Range("E8").Value = "  -0.02%  "
Set-TextValue "D9" "0.516"
$ws.Range("E9").Value = "  +0.05%  "
Set-TextValue "D10" "35.82"
$ws.Range("E10").Value = "  -1.50%  "
Set-TextValue "D11" "0.0792"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  -3.95%  "
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "2.665.40"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "2.292.92"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "43.007.69"
$ws.Range("E18").Value = "  -0.25%  "
Set-TextValue "D19" "13.11"
$ws.Range("E19").Value = "  +1.72%  "
$ws.Range("E20").Value = "  +0.26%  "
Set-TextValue "D21" "6.12"
$ws.Range("E21").Value = "  -1.09%  "
Set-TextValue "D22" "68.36"
$ws.Range("E22").Value = "  +0.35%  "
Set-TextValue "D23" "237.96"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  -0.54%  "
Set-TextValue "D30" "9.17"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("E31").Value = "  -13.59%  "
Set-TextValue "D32" "33.07"
$ws.Range("E32").Value = "  -5.32%  "
Set-TextValue "D33" "0.999"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("E34").Value = "  +2.52%  "
Set-TextValue "D35" "5.13"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -0.47%  "
Set-TextValue "D38" "0.0691"
$ws.Range("E38").Value = "  -0.85%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("E41").Value = "  +0.80%  "
Set-TextValue "D42" "2.75"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").Value = "2.009.69"
$ws.Range("E43").Value = "  +0.87%  "
Set-TextValue "D44" "0.0287"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("E45").Value = "  -7.59%  "
Set-TextValue "D46" "10.27"
$ws.Range("E46").Value = "  +1.20%  "
Set-TextValue "D47" "17.48"
$ws.Range("E47").Value = "  -1.36%  "
$ws.Range("E48").Value = "  -2.41%  "
Set-TextValue "D49" "54.48"
$ws.Range("D50").Value = "2.537.40"
$ws.Range("E50").Value = "  +0.36%  "
Set-TextValue "D51" "1.53"
$ws.Range("E51").Value = "  -1.36%  "
